$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts B:AN -> C:AO)
$ws.Columns("B:B").Insert()

# Set the new column header
$ws.Range("B1").Value = "Assignee ID"

# Match the new column's width to column A's width
$ws.Columns("B:B").ColumnWidth = $ws.Columns("A:A").ColumnWidth

# Refresh the autofilter so its range grows to include the new column
$ws.AutoFilterMode = $false
$ws.Range("A1:AO1").AutoFilter()

# Keep the workbook-level hidden _FilterDatabase name in sync with the new range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Data!_FilterDatabase") {
        $n.RefersTo = "=Data!`$A`$1:`$AO`$1"
    }
}

# Move the selection onto the newly inserted column's header cell
$ws.Range("B1").Select()
